$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "Oriolus_oriolus" (row 12)
$ws.Rows.Item(12).Delete()

# After the above deletion, "Sylvia_borin" (originally row 20) is now row 19
$ws.Rows.Item(19).Delete()
